# Apply crypto price/volume updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a plain decimal number would be
# auto-converted to a numeric type (losing the exact printed text, e.g.
# "216.30" -> 216.3) unless the cell is first forced to Text format.
$textForceCells = @(
    "D5",
    "D8",
    "D9",
    "D10",
    "D13",
    "D14",
    "D15",
    "D18",
    "D20",
    "D23",
    "D24",
    "D26",
    "D28",
    "D31",
    "D32",
    "D34",
    "D37",
    "D38",
    "D41",
    "D42",
    "D44",
    "D45",
    "D48",
    "D50",
    "D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.250.67'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.631.38'
$ws.Range('E3').Value = '  -1.20%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '216.30'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('E6').Value = '  +1.03%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.255'
$ws.Range('E8').Value = '  -0.45%  '
$ws.Range('D9').Value = '0.0625'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').Value = '20.29'
$ws.Range('E10').Value = '  +1.51%  '
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').Value = '1.629.44'
$ws.Range('E12').Value = '  -1.38%  '
$ws.Range('D13').Value = '4.12'
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('D14').Value = '0.544'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('D15').Value = '65.04'
$ws.Range('E15').Value = '  -3.98%  '
$ws.Range('D16').Value = '27.216.17'
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').Value = '217.04'
$ws.Range('E18').Value = '  -1.41%  '
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').Value = '6.93'
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('E22').Value = '  -4.14%  '
$ws.Range('D23').Value = '9.09'
$ws.Range('E23').Value = '  -1.60%  '
$ws.Range('D24').Value = '148.05'
$ws.Range('E24').Value = '  +0.78%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').Value = '7.26'
$ws.Range('E26').Value = '  -3.65%  '
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('D28').Value = '15.58'
$ws.Range('E28').Value = '  -1.52%  '
$ws.Range('E29').Value = '  -0.66%  '
$ws.Range('D31').Value = '3.38'
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('D32').Value = '2.99'
$ws.Range('E32').Value = '  -1.56%  '
$ws.Range('D33').Value = '1.318.21'
$ws.Range('E33').Value = '  +4.17%  '
$ws.Range('D34').Value = '1.56'
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('E35').Value = '  -0.31%  '
$ws.Range('E36').Value = '  -1.84%  '
$ws.Range('D37').Value = '0.541'
$ws.Range('E37').Value = '  -1.10%  '
$ws.Range('D38').Value = '0.847'
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('E40').Value = '  +2.05%  '
$ws.Range('D41').Value = '0.799'
$ws.Range('E41').Value = '  -1.29%  '
$ws.Range('D42').Value = '63.75'
$ws.Range('E42').Value = '  +2.56%  '
$ws.Range('D43').Value = '1.769.26'
$ws.Range('E43').Value = '  -1.28%  '
$ws.Range('D44').Value = '5.21'
$ws.Range('E44').Value = '  -4.48%  '
$ws.Range('D45').Value = '90.70'
$ws.Range('E45').Value = '  -1.22%  '
$ws.Range('E47').Value = '  +0.88%  '
$ws.Range('D48').Value = '0.811'
$ws.Range('E48').Value = '  +21.13%  '
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.0969'
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '7.54'
$ws.Range('E51').Value = '  -1.42%  '

Write-Output "Applied cryptos update"